# The commit deletes one gene row ("MSN1") from the "network" sheet,
# shifting all subsequent rows up by one (row 7 -> dimension A1:V22 becomes A1:V21).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network")
$ws.Activate()

# Select and delete the entire row 7 (MSN1), same as a user right-click > Delete
# on the row header - this shifts rows 8:22 up to 7:21 and updates the sheet
# dimension + shared-string usage count automatically.
$ws.Rows("7:7").Select() | Out-Null
$ws.Rows("7:7").Delete() | Out-Null
